$wb = $excel.ActiveWorkbook

# Update the Users sheet: replace "Nicole Bicho" with "Drew Koecher"
$usersSheet = $wb.Worksheets.Item("Users")
$usersSheet.Range("B4").Value = "Drew Koecher"

# Select B4 on the Users sheet and make Users the active sheet/tab
$usersSheet.Activate()
$usersSheet.Range("B4").Select()
